$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for Berenjena (Vega Modelo de Temuco),
# inserted as row 171; every existing record from the old row 171 onward
# shifts down by one row (old row 225 becomes row 226).
$ws.Rows.Item(171).Insert()

$ws.Cells.Item(171, 1).Value = 10
$ws.Cells.Item(171, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(171, 3).Value = "La Araucanía"
$ws.Cells.Item(171, 4).Value = 44588
$ws.Cells.Item(171, 5).Value = 9
$ws.Cells.Item(171, 6).Value = 100112001
$ws.Cells.Item(171, 7).Value = "Berenjena"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 65
$ws.Cells.Item(171, 11).Value = 10000
$ws.Cells.Item(171, 12).Value = 10000
$ws.Cells.Item(171, 13).Value = 10000
$ws.Cells.Item(171, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Región del Maule"
$ws.Cells.Item(171, 16).Value = 167
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"
